$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-5 (only columns C and E change)
$ws.Range("C2").Value = 23.42013439223523
$ws.Range("E2").Value = 0.5804712546134158

$ws.Range("C3").Value = 24.05579222437616
$ws.Range("E3").Value = 0.613742543171837

$ws.Range("C4").Value = 25.43849144375257
$ws.Range("E4").Value = 0.6842200997350938

$ws.Range("C5").Value = 29.09421076393417
$ws.Range("E5").Value = 0.8316405635363789

# New rows 6-17 with data (column A values are numeric-looking labels stored as text)
$newRows = @(
    @{ Row=6;  A="4";  B=23.61109032665684; C=22.7552839708626;  D=0.5630425821079695; E=0.5316978520688432 },
    @{ Row=7;  A="5";  B=23.56698259185531; C=23.31498834013616; D=0.560701678971915;  E=0.565062299647228 },
    @{ Row=8;  A="6";  B=23.61746968066501; C=24.7460004079343;  D=0.5623654053350033; E=0.6464005945857576 },
    @{ Row=9;  A="7";  B=23.67097614796583; C=28.35847660851668; D=0.5623341805466151; E=0.8031275800649054 },
    @{ Row=10; A="9";  B=23.58185794894158; C=28.29419717316968; D=0.5592424193153306; E=0.8042189327551276 },
    @{ Row=11; A="10"; B=23.53296335514261; C=24.70011184212721; D=0.5607680889382359; E=0.6485018772956437 },
    @{ Row=12; A="11"; B=23.53707245832457; C=23.336184782947;   D=0.5618308869854264; E=0.5707600925522629 },
    @{ Row=13; A="12"; B=23.62904466930573; C=22.87548817301593; D=0.563244054210409;  E=0.5375440827307707 },
    @{ Row=14; A="13"; B=23.60078928168511; C=29.09656297293069; D=0.5629599752400475; E=0.8334006987874809 },
    @{ Row=15; A="14"; B=23.59847117915199; C=25.4901777388197;  D=0.5619946590095188; E=0.6862067042922995 },
    @{ Row=16; A="15"; B=23.59772295605822; C=24.04720191284321; D=0.5586627486975914; E=0.6104305880609839 },
    @{ Row=17; A="16"; B=23.59866594249343; C=23.37781593296855; D=0.5600088516890556; E=0.5763392919183502 }
)

foreach ($r in $newRows) {
    $cellA = $ws.Cells.Item($r.Row, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $r.A

    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
}

# New row 18, only column A has a value (text label, not a number)
$cellA18 = $ws.Cells.Item(18, 1)
$cellA18.NumberFormat = "@"
$cellA18.Value = "730.1483306884766"
